$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New interleaved order of rows 16-23: alternate worker per period,
# periods ascending (1801..1804), Gustavo's salario basico unchanged,
# Antonio's salario basico updated to 781242.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "13357039"
$ws.Range("D16").Value = "GUSTAVO ORTIZ RODRIGUEZ"
$ws.Range("E16").Value = "1801"
$ws.Range("F16").Value = 120000
$ws.Range("G16").Value = 3000000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73162172"
$ws.Range("D17").Value = "ANTONIO JOSE ARRIETA MUSLASCO"
$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 30000
$ws.Range("G17").Value = 781242

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "13357039"
$ws.Range("D18").Value = "GUSTAVO ORTIZ RODRIGUEZ"
$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 120000
$ws.Range("G18").Value = 3000000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73162172"
$ws.Range("D19").Value = "ANTONIO JOSE ARRIETA MUSLASCO"
$ws.Range("E19").Value = "1802"
$ws.Range("F19").Value = 30000
$ws.Range("G19").Value = 781242

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "13357039"
$ws.Range("D20").Value = "GUSTAVO ORTIZ RODRIGUEZ"
$ws.Range("E20").Value = "1803"
$ws.Range("F20").Value = 120000
$ws.Range("G20").Value = 3000000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73162172"
$ws.Range("D21").Value = "ANTONIO JOSE ARRIETA MUSLASCO"
$ws.Range("E21").Value = "1803"
$ws.Range("F21").Value = 30000
$ws.Range("G21").Value = 781242

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "13357039"
$ws.Range("D22").Value = "GUSTAVO ORTIZ RODRIGUEZ"
$ws.Range("E22").Value = "1804"
$ws.Range("F22").Value = 120000
$ws.Range("G22").Value = 3000000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73162172"
$ws.Range("D23").Value = "ANTONIO JOSE ARRIETA MUSLASCO"
$ws.Range("E23").Value = "1804"
$ws.Range("F23").Value = 30000
$ws.Range("G23").Value = 781242

$wb.Save()
